$d = $word.ActiveDocument

$pairs = @(
    @("42-25=17", "89-12=77"),
    @("87-52=35", "53-37=16"),
    @("26+71=97", "33+4=37"),
    @("40+24=64", "17+18=35"),
    @("48-6=42", "73-27=46"),
    @("28+27=55", "33+16=49"),
    @("24+40=64", "39+28=67"),
    @("74-58=16", "5-3=2"),
    @("51-32=19", "97-36=61"),
    @("97-32=65", "95-91=4"),
    @("21+62=83", "9+29=38"),
    @("24+43=67", "88+1=89"),
    @("4+67=71", "68+9=77"),
    @("34+24=58", "5+81=86"),
    @("41-37=4", "35-25=10"),
    @("40-8=32", "72-29=43"),
    @("7+82=89", "97-25=72"),
    @("49-34=15", "97-0=97"),
    @("46-38=8", "62+19=81"),
    @("21+7=28", "54+17=71"),
    @("28+30=58", "47+3=50"),
    @("45+28=73", "55-17=38"),
    @("5+49=54", "69+30=99"),
    @("70-27=43", "29+51=80"),
    @("39+23=62", "4+53=57"),
    @("35+60=95", "77-51=26"),
    @("42+48=90", "20+54=74"),
    @("68-61=7", "5+18=23"),
    @("25+16=41", "26-21=5"),
    @("10+25=35", "89-46=43"),
    @("11+72=83", "22+29=51"),
    @("45-13=32", "68+20=88"),
    @("89-21=68", "67-51=16"),
    @("9+8=17", "66+12=78"),
    @("77-16=61", "36+38=74"),
    @("38+24=62", "6+56=62"),
    @("60-32=28", "20+46=66"),
    @("97-21=76", "87-20=67"),
    @("5+13=18", "96-78=18"),
    @("13+27=40", "51-43=8"),
    @("3+81=84", "71-57=14"),
    @("31+27=58", "58-54=4"),
    @("89-70=19", "82-81=1"),
    @("30+3=33", "25+32=57"),
    @("60-12=48", "13+40=53"),
    @("10-7=3", "77-13=64"),
    @("74-72=2", "74+4=78"),
    @("35+28=63", "13+26=39"),
    @("67+15=82", "33+45=78"),
    @("85-78=7", "85-34=51"),
    @("69-67=2", "23-16=7"),
    @("98-46=52", "30+54=84"),
    @("76-56=20", "65-25=40"),
    @("13+72=85", "20-10=10"),
    @("52+18=70", "33+62=95"),
    @("91-6=85", "27-8=19"),
    @("74-16=58", "50+19=69"),
    @("53+0=53", "4+2=6"),
    @("17-8=9", "87-84=3"),
    @("54-36=18", "54-48=6"),
    @("72-24=48", "21+25=46"),
    @("12+47=59", "89-36=53"),
    @("57-17=40", "74-11=63"),
    @("38+11=49", "80-42=38"),
    @("60-21=39", "92-21=71"),
    @("17+3=20", "63+9=72"),
    @("32+38=70", "64-41=23"),
    @("93-45=48", "40+33=73"),
    @("45-31=14", "74+13=87"),
    @("39-4=35", "31+55=86"),
    @("42-4=38", "61-20=41"),
    @("33+58=91", "69+6=75"),
    @("61-10=51", "14+3=17"),
    @("81-36=45", "29-26=3"),
    @("92-5=87", "63-38=25"),
    @("1+75=76", "92-19=73"),
    @("27+19=46", "72-0=72"),
    @("45+46=91", "0+67=67"),
    @("84-37=47", "25-6=19"),
    @("62+28=90", "98-11=87"),
    @("23+25=48", "5+18=23"),
    @("20+23=43", "86-3=83"),
    @("67+1=68", "64+3=67"),
    @("0+57=57", "39-31=8"),
    @("24+46=70", "55+22=77"),
    @("50-47=3", "95-21=74"),
    @("25+39=64", "23+51=74"),
    @("71-53=18", "49+21=70"),
    @("15+47=62", "59-31=28"),
    @("0+14=14", "27-26=1"),
    @("87-17=70", "35+7=42"),
    @("42+0=42", "81+8=89"),
    @("3+37=40", "89-38=51"),
    @("95-82=13", "68-19=49"),
    @("24+15=39", "74+5=79"),
    @("61-3=58", "56+8=64"),
    @("34+35=69", "66-56=10"),
    @("18+34=52", "47+36=83"),
    @("23+32=55", "91-76=15"),
    @("22-8=14", "40+26=66")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
